$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-04-01 Monday" "2024-04-02 Tuesday"

Replace-Text "282÷7=" "402÷7="
Replace-Text "697÷8=" "897÷2="
Replace-Text "560÷3=" "673÷4="
Replace-Text "168÷9=" "279÷4="
Replace-Text "181÷8=" "952÷6="
Replace-Text "257÷6=" "776÷4="
Replace-Text "541÷8=" "501÷7="
Replace-Text "177÷7=" "308÷2="
Replace-Text "731÷9=" "182÷8="
Replace-Text "125÷3=" "623÷7="
Replace-Text "268÷2=" "123÷6="
Replace-Text "722÷9=" "770÷2="
Replace-Text "474÷9=" "638÷7="
Replace-Text "532÷4=" "876÷7="
Replace-Text "522÷4=" "778÷9="
Replace-Text "583÷5=" "435÷4="
Replace-Text "656÷2=" "250÷2="
Replace-Text "530÷4=" "430÷9="
Replace-Text "476÷6=" "338÷9="
Replace-Text "829÷7=" "115÷6="
Replace-Text "102÷7=" "649÷4="
Replace-Text "275÷5=" "259÷7="
Replace-Text "829÷4=" "210÷4="
Replace-Text "993÷3=" "165÷4="
Replace-Text "723÷5=" "405÷7="
